$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns per the scraper refresh.
# Leading apostrophe forces text storage so values such as "309.30" and
# "0.45%" keep their exact original formatting (trailing zeros, literal
# percent sign) instead of being reinterpreted as numbers by Excel.
$ws.Range("D2").Value = "'309.30"
$ws.Range("E2").Value = "'0.45%"
$ws.Range("D3").Value = "'41.13"
$ws.Range("E3").Value = "'0.42%"
$ws.Range("D4").Value = "'5.232"
$ws.Range("E4").Value = "'2.50%"
$ws.Range("D5").Value = "'0.07689"
$ws.Range("E5").Value = "'0.81%"
$ws.Range("D6").Value = "'1.647"
$ws.Range("E6").Value = "'2.66%"
$ws.Range("E7").Value = "'1.34%"
$ws.Range("D8").Value = "'2.430"
$ws.Range("E8").Value = "'-1.65%"
$ws.Range("E9").Value = "'10.68%"
$ws.Range("D10").Value = "'0.1826"
$ws.Range("E10").Value = "'2.45%"
$ws.Range("D11").Value = "'0.09208"
$ws.Range("E11").Value = "'1.22%"
$ws.Range("D12").Value = "'0.04223"
$ws.Range("E12").Value = "'-0.05%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.001253"
$ws.Range("E14").Value = "'-0.22%"
$ws.Range("D15").Value = "'0.005736"
$ws.Range("E15").Value = "'-0.28%"
$ws.Range("E16").Value = "'1,903.02%"
$ws.Range("D18").Value = "'4.318"
$ws.Range("E18").Value = "'1.60%"
$ws.Range("D20").Value = "'7.394"
$ws.Range("E20").Value = "'11.38%"
$ws.Range("E22").Value = "'0.59%"
$ws.Range("D23").Value = "'0.04016"
$ws.Range("E23").Value = "'-1.47%"
$ws.Range("D24").Value = "'0.001266"
$ws.Range("E24").Value = "'1.69%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D38").Value = "'0.02553"
$ws.Range("E38").Value = "'7.44%"
$ws.Range("D39").Value = "'0.05339"
$ws.Range("E39").Value = "'3.29%"
$ws.Range("D40").Value = "'0.007837"
$ws.Range("E40").Value = "'0.61%"
$ws.Range("E41").Value = "'1.17%"
$ws.Range("D42").Value = "'0.006666"
$ws.Range("E42").Value = "'-5.52%"
$ws.Range("E43").Value = "'-4.71%"
$ws.Range("D44").Value = "'0.008022"
$ws.Range("E44").Value = "'0.99%"
$ws.Range("D45").Value = "'0.3074"
$ws.Range("E45").Value = "'-0.15%"
$ws.Range("D46").Value = "'0.00006727"
$ws.Range("E46").Value = "'-3.96%"
$ws.Range("E47").Value = "'-0.07%"
$ws.Range("D48").Value = "'0.1743"
$ws.Range("E48").Value = "'454.46%"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E51").Value = "'-0.07%"
